$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: classical-best-embeddings vs. classical-best-tfidf
$ws.Range("A2").Value = "classical-best-embed vs. classical-best-tfidf"
$ws.Range("C2").Value = 0.102
$ws.Range("D2").Value = 0.054
$ws.Range("I2").Value = 0.07199999999999999
$ws.Range("J2").Value = 0.07099999999999999

# Row 3: BERT-base vs. classical-best-tfidf
$ws.Range("C3").Value = 0.091
$ws.Range("D3").Value = 0.158
$ws.Range("E3").Value = 0.154
$ws.Range("F3").Value = 0.151
$ws.Range("G3").Value = 0.157
$ws.Range("H3").Value = 0.183
$ws.Range("I3").Value = 0.138
$ws.Range("J3").Value = 0.149

# Row 4: BERT-base vs. classical-best-embeddings
$ws.Range("A4").Value = "BERT-base vs. classical-best-embed"
$ws.Range("C4").Value = -0.011
$ws.Range("D4").Value = 0.104
$ws.Range("E4").Value = 0.08500000000000001
$ws.Range("F4").Value = 0.089
$ws.Range("G4").Value = 0.104
$ws.Range("H4").Value = 0.095
$ws.Range("I4").Value = 0.067
$ws.Range("J4").Value = 0.078

# Row 5: BERT-base-nli vs. classical-best-tfidf
$ws.Range("B5").Value = 0.367
$ws.Range("C5").Value = 0.258
$ws.Range("D5").Value = 0.223
$ws.Range("E5").Value = 0.218
$ws.Range("F5").Value = 0.207
$ws.Range("G5").Value = 0.195
$ws.Range("H5").Value = 0.212
$ws.Range("I5").Value = 0.226
$ws.Range("J5").Value = 0.219

# Row 6: BERT-base-nli vs. classical-best-embeddings
$ws.Range("A6").Value = "BERT-base-nli vs. classical-best-embed"
$ws.Range("B6").Value = 0.367
$ws.Range("C6").Value = 0.156
$ws.Range("D6").Value = 0.169
$ws.Range("E6").Value = 0.149
$ws.Range("F6").Value = 0.145
$ws.Range("G6").Value = 0.142
$ws.Range("H6").Value = 0.124
$ws.Range("I6").Value = 0.155
$ws.Range("J6").Value = 0.148

# Row 7: BERT-base-nli vs. BERT-base
$ws.Range("B7").Value = 0.367
$ws.Range("C7").Value = 0.167
$ws.Range("D7").Value = 0.065
$ws.Range("E7").Value = 0.064
$ws.Range("F7").Value = 0.056
$ws.Range("G7").Value = 0.038
$ws.Range("H7").Value = 0.029
$ws.Range("I7").Value = 0.08799999999999999
